$wb = $excel.ActiveWorkbook

# ===========================================================================
# 1) Update the existing "ODI Batting" sheet:
#    - rename MATCH_CARD_LINK header to MATCH_CODE
#    - replace the full scorecard URL with just the numeric match code (text)
#    - clear the handful of truly-empty INNING_NUMBER cells (did-not-bat rows)
# ===========================================================================
$wsBatting = $wb.Worksheets.Item("ODI Batting")

$wsBatting.Range("D1").Value = "MATCH_CODE"

$wsBatting.Range("D2:D24").NumberFormat = "@"
$battingCodes = @{
    2  = "4154"; 3  = "4159"; 4  = "4161"; 5  = "4164"; 6  = "4180";
    7  = "4181"; 8  = "4220"; 9  = "4221"; 10 = "4228"; 11 = "4229";
    12 = "4230"; 13 = "4362"; 14 = "4379"; 15 = "4385"; 16 = "4387";
    17 = "4388"; 18 = "4391"; 19 = "4413"; 20 = "4414"; 21 = "4592";
    22 = "4611"; 23 = "4616"; 24 = "4624"
}
foreach ($row in $battingCodes.Keys) {
    $wsBatting.Cells.Item($row, 4).Value = $battingCodes[$row]
}

$emptyInningRows = @(3, 7, 14, 15, 17, 18)
foreach ($row in $emptyInningRows) {
    $wsBatting.Cells.Item($row, 2).ClearContents()
}

# ===========================================================================
# 2) Update the existing "ODI Bowling" sheet:
#    - rename MATCH_CARD_LINK header to MATCH_CODE
#    - replace the full scorecard URL with just the numeric match code (text)
# ===========================================================================
$wsBowling = $wb.Worksheets.Item("ODI Bowling")

$wsBowling.Range("B1").Value = "MATCH_CODE"

$wsBowling.Range("B2:B23").NumberFormat = "@"
$bowlingCodes = @{
    2  = "4154"; 3  = "4159"; 4  = "4161"; 5  = "4164"; 6  = "4180";
    7  = "4181"; 8  = "4220"; 9  = "4221"; 10 = "4228"; 11 = "4229";
    12 = "4230"; 13 = "4362"; 14 = "4379"; 15 = "4385"; 16 = "4387";
    17 = "4388"; 18 = "4391"; 19 = "4413"; 20 = "4414"; 21 = "4592";
    22 = "4611"; 23 = "4624"
}
foreach ($row in $bowlingCodes.Keys) {
    $wsBowling.Cells.Item($row, 2).Value = $bowlingCodes[$row]
}

# ===========================================================================
# 3) Insert a new "Player Info" sheet before "ODI Batting"
# ===========================================================================
$wsBattingAnchor = $wb.Worksheets.Item("ODI Batting")
$wsPlayerInfo = $wb.Worksheets.Add($wsBattingAnchor)
$wsPlayerInfo.Name = "Player Info"

$wsPlayerInfo.Range("A1").Value = "ID"
$wsPlayerInfo.Range("B1").Value = "NAME"
$wsPlayerInfo.Range("C1").Value = "BATTING_HAND"
$wsPlayerInfo.Range("D1").Value = "BOWL_STYLE"

$wsPlayerInfo.Range("A1:D1").Font.Bold = $true
$wsPlayerInfo.Range("A1:D1").HorizontalAlignment = -4108
$wsPlayerInfo.Range("A1:D1").VerticalAlignment = -4160
$wsPlayerInfo.Range("A1:D1").Borders.LineStyle = 1

$wsPlayerInfo.Range("A2").NumberFormat = "@"
$wsPlayerInfo.Range("A2").Value = "4742"
$wsPlayerInfo.Range("B2").Value = "Keemo Mandela Angus Paul"
$wsPlayerInfo.Range("C2").Value = "Right Handed"
$wsPlayerInfo.Range("D2").Value = "Right Arm Fast Medium"

# ===========================================================================
# 4) Append a new "ODI Batting Extra" sheet after "ODI Bowling"
# ===========================================================================
$wsBowlingAnchor = $wb.Worksheets.Item("ODI Bowling")
$wsExtra = $wb.Worksheets.Add($null, $wsBowlingAnchor)
$wsExtra.Name = "ODI Batting Extra"

$wsExtra.Range("A1").Value = "MATCH_CODE"
$wsExtra.Range("B1").Value = "BATTING_POSITION"
$wsExtra.Range("C1").Value = "NUM_4"
$wsExtra.Range("D1").Value = "NUM_6"
$wsExtra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$wsExtra.Range("F1").Value = "MAN_OF_MATCH"

$wsExtra.Range("A1:F1").Font.Bold = $true
$wsExtra.Range("A1:F1").HorizontalAlignment = -4108
$wsExtra.Range("A1:F1").VerticalAlignment = -4160
$wsExtra.Range("A1:F1").Borders.LineStyle = 1

# MATCH_CODE, NUM_4, NUM_6 and PERCENT_RUNS_OF_TOTAL are text columns;
# BATTING_POSITION stays numeric.
$wsExtra.Range("A2:A21").NumberFormat = "@"
$wsExtra.Range("C2:E21").NumberFormat = "@"

$extraData = @(
    @{ Row = 2;  Code = "4164"; Pos = 10; Num4 = "0"; Num6 = "0"; Pct = $null;    Mom = "NO" },
    @{ Row = 3;  Code = "4180"; Pos = 9;  Num4 = "1"; Num6 = "0"; Pct = "1.48%";  Mom = "NO" },
    @{ Row = 4;  Code = "4181"; Pos = 9;  Num4 = $null; Num6 = $null; Pct = $null; Mom = "NO" },
    @{ Row = 5;  Code = "4220"; Pos = 10; Num4 = "1"; Num6 = "2"; Pct = "12.42%"; Mom = "NO" },
    @{ Row = 6;  Code = "4221"; Pos = 8;  Num4 = "0"; Num6 = "0"; Pct = "4.81%";  Mom = "NO" },
    @{ Row = 7;  Code = "4228"; Pos = 8;  Num4 = "1"; Num6 = "2"; Pct = "18.46%"; Mom = "NO" },
    @{ Row = 8;  Code = "4229"; Pos = 8;  Num4 = "0"; Num6 = "0"; Pct = "7.03%";  Mom = "NO" },
    @{ Row = 9;  Code = "4230"; Pos = 9;  Num4 = "0"; Num6 = "0"; Pct = "6.06%";  Mom = "NO" },
    @{ Row = 10; Code = "4362"; Pos = $null; Num4 = $null; Num6 = $null; Pct = $null; Mom = "NO" },
    @{ Row = 11; Code = "4379"; Pos = 8;  Num4 = $null; Num6 = $null; Pct = $null; Mom = "NO" },
    @{ Row = 12; Code = "4385"; Pos = $null; Num4 = $null; Num6 = $null; Pct = $null; Mom = "NO" },
    @{ Row = 13; Code = "4387"; Pos = 8;  Num4 = "4"; Num6 = "3"; Pct = "16.43%"; Mom = "NO" },
    @{ Row = 14; Code = "4388"; Pos = 8;  Num4 = $null; Num6 = $null; Pct = $null; Mom = "NO" },
    @{ Row = 15; Code = "4391"; Pos = 8;  Num4 = $null; Num6 = $null; Pct = $null; Mom = "NO" },
    @{ Row = 16; Code = "4413"; Pos = 8;  Num4 = "4"; Num6 = "1"; Pct = "11.07%"; Mom = "NO" },
    @{ Row = 17; Code = "4414"; Pos = $null; Num4 = $null; Num6 = $null; Pct = $null; Mom = "NO" },
    @{ Row = 18; Code = "4592"; Pos = 8;  Num4 = "2"; Num6 = "1"; Pct = "9.72%";  Mom = "NO" },
    @{ Row = 19; Code = "4611"; Pos = 8;  Num4 = "4"; Num6 = "0"; Pct = "23.15%"; Mom = "NO" },
    @{ Row = 20; Code = "4616"; Pos = 7;  Num4 = "1"; Num6 = "0"; Pct = "3.37%";  Mom = "NO" },
    @{ Row = 21; Code = "4624"; Pos = 9;  Num4 = "0"; Num6 = "0"; Pct = $null;    Mom = "NO" }
)

foreach ($item in $extraData) {
    $r = $item.Row
    $wsExtra.Cells.Item($r, 1).Value = $item.Code
    if ($null -ne $item.Pos) { $wsExtra.Cells.Item($r, 2).Value = $item.Pos }
    if ($null -ne $item.Num4) { $wsExtra.Cells.Item($r, 3).Value = $item.Num4 }
    if ($null -ne $item.Num6) { $wsExtra.Cells.Item($r, 4).Value = $item.Num6 }
    if ($null -ne $item.Pct) { $wsExtra.Cells.Item($r, 5).Value = $item.Pct }
    $wsExtra.Cells.Item($r, 6).Value = $item.Mom
}
